# Append the new devlog entry (row 8) to Sheet1, mirroring rows 2-7:
#   A: date (styled like the other date cells)
#   B: task description (new shared string)
#   C: hours worked

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data - 26 Nov 2024 (serial 45622), matching the date style
# already used by A2:A7.
$ws.Range("A8").Value = 45622
$ws.Range("A8").NumberFormat = $ws.Range("A7").NumberFormat

$ws.Range("B8").Value = "Some more front end for dialogue and a beginning on the inventory."

$ws.Range("C8").Value = 3

# Excel moves the active selection to the newly-edited cell.
[void]$ws.Range("B8").Select()

Write-Output "Added row 8 devlog entry"
